$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The agency's roster (rows 2-6) is consolidated onto a single contact:
# every row's Name becomes "edwin camacho" and every row's phone number
# becomes the text value "300 6120261" (previously a numeric phone number
# shared by placeholder names).
$ws.Range("A2").Value = "edwin camacho"
$ws.Range("B2").Value = "300 6120261"

$ws.Range("A3").Value = "edwin camacho"
$ws.Range("B3").Value = "300 6120261"

$ws.Range("A4").Value = "edwin camacho"
$ws.Range("B4").Value = "300 6120261"

$ws.Range("A5").Value = "edwin camacho"
$ws.Range("B5").Value = "300 6120261"

$ws.Range("A6").Value = "edwin camacho"
$ws.Range("B6").Value = "300 6120261"

# Reflect the new selection shown in the sheet view (user had selected
# B4:B6 after editing the phone numbers).
$ws.Range("B4:B6").Select()
